# Weekly update: a new week's price record is inserted at the top of the
# data table (row 7), pushing all existing data rows down by one. The
# previously-last row (old row 138) becomes the new last row (139).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 7 (first data row, right below the header).
# This shifts rows 7..138 down to 8..139 and extends the used range.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row 7 with this week's record.
$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C7").Value = "Ñuble"
$ws.Range("D7").Value = 44496
$ws.Range("E7").Value = 16
$ws.Range("F7").Value = 100112006
$ws.Range("G7").Value = "Repollo"
$ws.Range("H7").Value = "Crespo record"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 600
$ws.Range("L7").Value = 700
$ws.Range("M7").Value = 650
$ws.Range("N7").Value = "$/unidad"
$ws.Range("O7").Value = "Región del Maule"
$ws.Range("P7").Value = 650
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = "Hortaliza"
